# CDM_Servitudes.xlsx — header row rework
#
# The sheet's header row (row 1) is reordered/relabeled:
#   old: A=LandId B=ServitudeId C=Name       D=ValidFrom E=ValidUntil F=GUID
#   new: A=ServitudeId B=Type C=ValidFrom D=ValidUntil E=LandId F=Guid
# ("Name" is dropped, "Type" is a new column, "GUID" is renamed "Guid").
# The header row is then bolded and selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "ServitudeId"
$ws.Cells.Item(1, 2).Value = "Type"
$ws.Cells.Item(1, 3).Value = "ValidFrom"
$ws.Cells.Item(1, 4).Value = "ValidUntil"
$ws.Cells.Item(1, 5).Value = "LandId"
$ws.Cells.Item(1, 6).Value = "Guid"

$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true

[void]$headerRange.Select()
